# Fix the spelling of the document title: "Navagación" -> "Navegación"
# and move the "_GoBack" bookmark (Word's "last edit location" marker)
# from the end of the document to the point of this correction, exactly
# as real Word would do after typing the fix in place.

$d = $word.ActiveDocument

# The existing "_GoBack" bookmark currently sits at the end of the
# document (after the last bullet point). Remove it - it will be
# re-created at the location of the correction below.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# The title is the very first paragraph: "Navagación" (chars 0-9).
#   N(0) a(1) v(2) a(3) g(4) a(5) c(6) i(7) o(8) n(9)
# Only the "a" at offset 3 needs to become "e".
$rTypo = $d.Range(3, 4)
$rTypo.Text = "e"

# Drop a bookmark right after the corrected "e" (offset 4), between it
# and "gación" - this reproduces the "_GoBack" bookmark Word leaves at
# the last edited spot.
$rGoBack = $d.Range(4, 4)
$d.Bookmarks.Add("_GoBack", $rGoBack)

# Also split "Nav" from "e" into separate runs (as Word's editor would
# leave behind from the incremental keystrokes), using a temporary
# bookmark to force the run boundary, then discarding it.
$rSplit = $d.Range(3, 3)
$d.Bookmarks.Add("zzzTempSplit", $rSplit)
$d.Bookmarks.Item("zzzTempSplit").Delete()
